$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '23.314.97'
Set-TextCell $ws 'E2' '  -0.42%  '

Set-TextCell $ws 'D3' '1.625.52'
Set-TextCell $ws 'E3' '  -0.66%  '

Set-TextCell $ws 'E4' '  +0.03%  '

Set-TextCell $ws 'D5' '1.002'
Set-TextCell $ws 'E5' '  +0.01%  '

Set-TextCell $ws 'D6' '303.21'
Set-TextCell $ws 'E6' '  -0.59%  '

Set-TextCell $ws 'E7' '  +0.16%  '

Set-TextCell $ws 'D8' '0.3623'
Set-TextCell $ws 'E8' '  +0.18%  '

Set-TextCell $ws 'D9' '51.31'
Set-TextCell $ws 'E9' '  -1.13%  '

Set-TextCell $ws 'D10' '0.08151'
Set-TextCell $ws 'E10' '  +0.42%  '

Set-TextCell $ws 'D11' '1.225'
Set-TextCell $ws 'E11' '  -2.18%  '

Set-TextCell $ws 'D12' '1.001'
Set-TextCell $ws 'E12' '  -0.15%  '

Set-TextCell $ws 'D13' '22.25'
Set-TextCell $ws 'E13' '  -2.45%  '

Set-TextCell $ws 'D14' '6.472'
Set-TextCell $ws 'E14' '  -1.93%  '

Set-TextCell $ws 'B15' 'Chainlink'
Set-TextCell $ws 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell $ws 'D15' '7.301'
Set-TextCell $ws 'E15' '  +0.35%  '

Set-TextCell $ws 'B16' 'ShibaInu'
Set-TextCell $ws 'C16' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D16' '0.00001237'
Set-TextCell $ws 'E16' '  -2.20%  '

Set-TextCell $ws 'D17' '1.619.06'
Set-TextCell $ws 'E17' '  -1.11%  '

Set-TextCell $ws 'D18' '93.82'
Set-TextCell $ws 'E18' '  -0.46%  '

Set-TextCell $ws 'D19' '0.06956'
Set-TextCell $ws 'E19' '  +0.75%  '

Set-TextCell $ws 'D20' '17.50'
Set-TextCell $ws 'E20' '  -3.32%  '

Set-TextCell $ws 'D21' '6.549'
Set-TextCell $ws 'E21' '  +0.64%  '

Set-TextCell $ws 'D22' '1.002'
Set-TextCell $ws 'E22' '  +0.11%  '

Set-TextCell $ws 'D23' '12.52'
Set-TextCell $ws 'E23' '  -1.55%  '

Set-TextCell $ws 'D24' '23.351.63'
Set-TextCell $ws 'E24' '  -0.27%  '

Set-TextCell $ws 'E25' '  +2.56%  '

Set-TextCell $ws 'D26' '2.466'
Set-TextCell $ws 'E26' '  +1.85%  '

Set-TextCell $ws 'D27' '21.20'
Set-TextCell $ws 'E27' '  +0.19%  '

Set-TextCell $ws 'D28' '150.16'
Set-TextCell $ws 'E28' '  -0.98%  '

Set-TextCell $ws 'D29' '5.242'
Set-TextCell $ws 'E29' '  -1.42%  '

Set-TextCell $ws 'D30' '132.58'
Set-TextCell $ws 'E30' '  -2.44%  '

Set-TextCell $ws 'D31' '1.801.40'
Set-TextCell $ws 'E31' '  -0.93%  '

Set-TextCell $ws 'B32' 'Filecoin'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D32' '6.734'
Set-TextCell $ws 'E32' '  +0.20%  '

Set-TextCell $ws 'B33' 'WEMIXTOKEN'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws 'D33' '2.178'
Set-TextCell $ws 'E33' '  -4.76%  '

Set-TextCell $ws 'D34' '1.039'
Set-TextCell $ws 'E34' '  +9.31%  '

Set-TextCell $ws 'D35' '10.77'
Set-TextCell $ws 'E35' '  +5.23%  '

Set-TextCell $ws 'D36' '0.02753'
Set-TextCell $ws 'E36' '  -2.38%  '

Set-TextCell $ws 'D37' '0.2494'
Set-TextCell $ws 'E37' '  -0.73%  '

Set-TextCell $ws 'D38' '0.08779'
Set-TextCell $ws 'E38' '  -0.07%  '

Set-TextCell $ws 'D39' '0.07112'
Set-TextCell $ws 'E39' '  -1.90%  '

Set-TextCell $ws 'D40' '5.979'
Set-TextCell $ws 'E40' '  -1.48%  '

Set-TextCell $ws 'D41' '0.6977'
Set-TextCell $ws 'E41' '  -0.78%  '

Set-TextCell $ws 'D42' '1.337'
Set-TextCell $ws 'E42' '  -2.41%  '

Set-TextCell $ws 'D43' '15.94'
Set-TextCell $ws 'E43' '  -1.26%  '

Set-TextCell $ws 'D44' '12.07'
Set-TextCell $ws 'E44' '  -2.62%  '

Set-TextCell $ws 'D45' '0.6476'
Set-TextCell $ws 'E45' '  -0.44%  '

Set-TextCell $ws 'D46' '1.001'
Set-TextCell $ws 'E46' '  +0.01%  '

Set-TextCell $ws 'B47' 'PancakeSwap'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws 'D47' '3.965'
Set-TextCell $ws 'E47' '  -1.07%  '

Set-TextCell $ws 'B48' 'NEARProtocol'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws 'D48' '2.264'
Set-TextCell $ws 'E48' '  -2.36%  '

Set-TextCell $ws 'D49' '0.07967'
Set-TextCell $ws 'E49' '  -0.05%  '

Set-TextCell $ws 'D50' '1.184'
Set-TextCell $ws 'E50' '  -1.08%  '

Set-TextCell $ws 'D51' '125.39'
Set-TextCell $ws 'E51' '  -2.12%  '
